$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for rows 2..23. Values below reflect the post-edit state of the
# workbook (rows reshuffled / updated per commit "Fruta / hortaliza, semanal").
$rows = @{
  2  = @{ D = 44607; J = 900;  K = 1300; L = 1400; M = 1350; P = 1350 }
  3  = @{ D = 44656; J = 1000; K = 900;  L = 1000; M = 950;  P = 950  }
  4  = @{ D = 44687; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 }
  5  = @{ D = 44455; J = 1100; K = 900;  L = 1000; M = 950;  P = 950  }
  6  = @{ D = 44550; J = 1300; K = 1000; L = 1200; M = 1100; P = 1100 }
  7  = @{ D = 44175; J = 1600; K = 1000; L = 1200; M = 1100; P = 1100 }
  8  = @{ D = 44673; J = 900;  K = 1300; L = 1400; M = 1350; P = 1350 }
  9  = @{ D = 44883; J = 800;  K = 550;  L = 600;  M = 575;  P = 575  }
  10 = @{ D = 44784; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 }
  11 = @{ D = 44638; J = 1000; K = 900;  L = 950;  M = 925;  P = 925  }
  12 = @{ D = 44243; J = 1200; K = 1200; L = 1300; M = 1250; P = 1250 }
  13 = @{ D = 44229; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
  14 = @{ D = 44449; J = 1300; K = 900;  L = 950;  M = 925;  P = 925  }
  15 = @{ D = 44341; J = 1300; K = 900;  L = 1000; M = 950;  P = 950  }
  16 = @{ D = 44407; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 }
  17 = @{ D = 44291; J = 1000; K = 1000; L = 1200; M = 1100; P = 1100 }
  18 = @{ D = 44649; J = 600;  K = 900;  L = 1000; M = 950;  P = 950  }
  19 = @{ D = 44476; J = 900;  K = 700;  L = 800;  M = 750;  P = 750  }
  20 = @{ D = 44453; J = 1000; K = 800;  L = 900;  M = 850;  P = 850  }
  21 = @{ D = 44442; J = 1250; K = 850;  L = 900;  M = 875;  P = 875  }
  22 = @{ D = 44284; J = 1500; K = 800;  L = 850;  M = 825;  P = 825  }
  23 = @{ D = 44484; J = 900;  K = 750;  L = 800;  M = 775;  P = 775  }
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Range("D$r").Value = $vals.D
  $ws.Range("J$r").Value = $vals.J
  $ws.Range("K$r").Value = $vals.K
  $ws.Range("L$r").Value = $vals.L
  $ws.Range("M$r").Value = $vals.M
  $ws.Range("P$r").Value = $vals.P
}
